# Update cryptocurrency values on the "Cripto" sheet with new (falling) prices.
# Values are stored/kept as plain text (Brazilian-style "1.234,56" formatting),
# matching how the workbook already represents them (shared strings, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cripto")

# Bitcoin: Valor Atual / Valor Anterior
$ws.Range("B2").Value = "28.497,8"
$ws.Range("C2").Value = "28.446,1"

# Ethereum: Valor Atual / Valor Anterior
$ws.Range("B3").Value = "1.537,2"
$ws.Range("C3").Value = "1.525,38"

# Dogecoin: Valor Atual / Valor Anterior
# These look like plain numbers once the comma is stripped (two groups of
# three digits, e.g. "069827"), so a bare assignment would be auto-converted
# to a number by the locale-aware parser. Force text storage by temporarily
# switching the cell to a Text number format, then restore the original
# ("General"/unset) format so the cell's style stays exactly as it was.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0,069827"
$ws.Range("B4").NumberFormat = ""

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0,079035"
$ws.Range("C4").NumberFormat = ""
